# Update attendance ("想去人数", column F) and minimum ticket price
# ("最低票价", column G) figures to the latest scraped values.
#
# Sheet "展览" and "全部类型" share the same exhibition rows (since
# "全部类型" is a merge of "展览" + "演出"), so every exhibition row
# touched here is updated on both sheets; the single "演出" row lives on
# "演出" and on "全部类型".

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" ----------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")

$wsExpo.Range("F2").Value = 629
$wsExpo.Range("G2").Value = 55

$wsExpo.Range("F4").Value = 47

$wsExpo.Range("F5").Value = 4586

$wsExpo.Range("F6").Value = 1829

$wsExpo.Range("F14").Value = 520

$wsExpo.Range("F21").Value = 1567

$wsExpo.Range("F31").Value = 3661

$wsExpo.Range("F32").Value = 753

$wsExpo.Range("F34").Value = 345

$wsExpo.Range("F35").Value = 55

$wsExpo.Range("F36").Value = 1754

# ---- Sheet "演出" ----------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")

$wsShow.Range("F2").Value = 21

# ---- Sheet "全部类型" (merged view of 展览 + 演出) -------------------
$wsAll = $wb.Worksheets.Item("全部类型")

$wsAll.Range("F2").Value = 629
$wsAll.Range("G2").Value = 55

$wsAll.Range("F4").Value = 47

$wsAll.Range("F5").Value = 4586

$wsAll.Range("F6").Value = 1829

$wsAll.Range("F14").Value = 520

$wsAll.Range("F16").Value = 21

$wsAll.Range("F22").Value = 1567

$wsAll.Range("F32").Value = 3661

$wsAll.Range("F34").Value = 753

$wsAll.Range("F36").Value = 345

$wsAll.Range("F37").Value = 55

$wsAll.Range("F38").Value = 1754
